# Generate Report for Archive
#
# 1. Replace the status text "Ready for handoff" with "In Translation"
#    wherever it occurs (Overview / zh-cn / de-de sheets).
# 2. Narrow the "Status"-related columns from ~17.22 to ~13.41 (character)
#    units wide on all three sheets.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value()
        if ("Ready for handoff" -eq $val) {
            $cell.Value = "In Translation"
        }
    }
}

# Narrow columns (values chosen so the resulting column width lands on the
# nearest achievable pixel-quantized width to the target 13.41 char units).
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
